$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "土地" (land) — first sheet. Previously only had the header row
# (B1:Q1). Add the first data row (row 2) describing a newly-disclosed
# land parcel.
# ---------------------------------------------------------------------------
$wsLand = $wb.Worksheets.Item(1)

$wsLand.Range("A2").Value = 14
$wsLand.Range("B2").Value = "基隆市安樂區大武崙段内寮小段05090010地號"
$wsLand.Range("C2").Value = 198
$wsLand.Range("D2").Value = "全部"
$wsLand.Range("E2").Value = "謝國樑"
$wsLand.Range("F2").Value = "91年12月26日"
$wsLand.Range("G2").Value = "受贈"
$wsLand.Range("H2").Value = 336600
$wsLand.Range("I2").Value = "land"
$wsLand.Range("J2").Value = "normal"
$wsLand.Range("K2").Value = "2011-11-17"
$wsLand.Range("L2").Value = "謝國樑"
$wsLand.Range("M2").Value = 1387
$wsLand.Range("N2").Value = "tmpbbad1"
$wsLand.Range("O2").Value = 14
$wsLand.Range("P2").Value = 1
$wsLand.Range("Q2").Value = 198

# Column A on the data rows carries the same bold/bordered/centered look as
# the header row (matches the existing workbook's quirky style reuse).
$wsLand.Range("A2").Font.Bold = $true
$wsLand.Range("A2").Borders.LineStyle = 1
$wsLand.Range("A2").HorizontalAlignment = -4108
$wsLand.Range("A2").VerticalAlignment = -4160

# ---------------------------------------------------------------------------
# Sheet "存款" (deposits) — insert a brand-new first data row (row 2) for a
# previously-missing bank account, pushing the existing rows down by one.
# ---------------------------------------------------------------------------
$wsDeposit = $wb.Worksheets.Item(2)

$wsDeposit.Rows.Item(2).Insert()
$wsDeposit.Range("B2:F2").ClearFormats()

$wsDeposit.Range("A2").Value = 43
$wsDeposit.Range("B2").Value = "基隆二信營業部"
$wsDeposit.Range("C2").Value = "活期儲蓄存款"
$wsDeposit.Range("D2").Value = "新臺幣"
$wsDeposit.Range("E2").Value = "謝國樑"
$wsDeposit.Range("F2").Value = 1

$wsDeposit.Range("A2").Font.Bold = $true
$wsDeposit.Range("A2").Borders.LineStyle = 1
$wsDeposit.Range("A2").HorizontalAlignment = -4108
$wsDeposit.Range("A2").VerticalAlignment = -4160

# ---------------------------------------------------------------------------
# Sheet "股票" (stocks) — insert a brand-new first data row (row 2) for a
# previously-missing stock holding, pushing the existing rows down by one.
# ---------------------------------------------------------------------------
$wsStock = $wb.Worksheets.Item(3)

$wsStock.Rows.Item(2).Insert()
$wsStock.Range("B2:N2").ClearFormats()

$wsStock.Range("A2").Value = 62
$wsStock.Range("B2").Value = "大魯閣纖维"
$wsStock.Range("C2").Value = "謝國樑"
$wsStock.Range("D2").Value = 22355
$wsStock.Range("E2").Value = 10
$wsStock.Range("F2").Value = "新臺幣"
$wsStock.Range("G2").Value = 223550
$wsStock.Range("H2").Value = "stock"
$wsStock.Range("I2").Value = "normal"
$wsStock.Range("J2").Value = "2011-11-17"
$wsStock.Range("K2").Value = "謝國樑"
$wsStock.Range("L2").Value = 1387
$wsStock.Range("M2").Value = "tmpbbad1"
$wsStock.Range("N2").Value = 62

$wsStock.Range("A2").Font.Bold = $true
$wsStock.Range("A2").Borders.LineStyle = 1
$wsStock.Range("A2").HorizontalAlignment = -4108
$wsStock.Range("A2").VerticalAlignment = -4160

Write-Output "property building done"
